# Update cryptocurrency price/volume figures (and the Bittensor/Mantle row swap)
# as scraped by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '54.778.59'
$ws.Range('E2').Value = '  +0.91%  '
# Row 3
$ws.Range('D3').Value = '2.304.39'
$ws.Range('E3').Value = '  +0.61%  '
# Row 4
$ws.Range('E4').Value = '  +0.06%  '
# Row 5
$ws.Range('D5').Value = '''498.25'
$ws.Range('E5').Value = '  +0.76%  '
# Row 6
$ws.Range('D6').Value = '''129.14'
$ws.Range('E6').Value = '  +0.69%  '
# Row 7
$ws.Range('E7').Value = '  +0.21%  '
# Row 8
$ws.Range('E8').Value = '  +0.83%  '
# Row 9
$ws.Range('D9').Value = '2.305.99'
$ws.Range('E9').Value = '  +0.31%  '
# Row 10
$ws.Range('D10').Value = '''0.0955'
$ws.Range('E10').Value = '  +1.73%  '
# Row 11
$ws.Range('E11').Value = '  +2.33%  '
# Row 12
$ws.Range('E12').Value = '  +2.23%  '
# Row 13
$ws.Range('E13').Value = '  -1.88%  '
# Row 14
$ws.Range('D14').Value = '2.710.64'
$ws.Range('E14').Value = '  +0.59%  '
# Row 15
$ws.Range('D15').Value = '''21.92'
$ws.Range('E15').Value = '  +2.31%  '
# Row 16
$ws.Range('D16').Value = '54.721.63'
$ws.Range('E16').Value = '  +0.85%  '
# Row 17
$ws.Range('E17').Value = '  +0.68%  '
# Row 18
$ws.Range('D18').Value = '2.277.25'
$ws.Range('E18').Value = '  -0.40%  '
# Row 19
$ws.Range('D19').Value = '''10.12'
$ws.Range('E19').Value = '  +4.44%  '
# Row 20
$ws.Range('E20').Value = '  +2.30%  '
# Row 21
$ws.Range('D21').Value = '''308.16'
$ws.Range('E21').Value = '  +1.46%  '
# Row 22
$ws.Range('E22').Value = '  +4.60%  '
# Row 23
$ws.Range('E23').Value = '  -0.01%  '
# Row 24
$ws.Range('E24').Value = '  -1.14%  '
# Row 25
$ws.Range('D25').Value = '''62.95'
$ws.Range('E25').Value = '  -1.90%  '
# Row 26
$ws.Range('E26').Value = '  +0.16%  '
# Row 27
$ws.Range('E27').Value = '  +5.83%  '
# Row 28
$ws.Range('E28').Value = '  +2.16%  '
# Row 29
$ws.Range('D29').Value = '2.406.74'
$ws.Range('E29').Value = '  +0.20%  '
# Row 30
$ws.Range('D30').Value = '''7.18'
$ws.Range('E30').Value = '  +0.76%  '
# Row 31
$ws.Range('D31').Value = '''169.70'
$ws.Range('E31').Value = '  +0.15%  '
# Row 32
$ws.Range('D32').Value = '0.0₃0701'
$ws.Range('E32').Value = '  -0.01%  '
# Row 33
$ws.Range('D33').Value = '''1.61'
$ws.Range('E33').Value = '  -0.27%  '
# Row 34
$ws.Range('E34').Value = '  +2.27%  '
# Row 36
$ws.Range('D36').Value = '''1.08'
$ws.Range('E36').Value = '  +0.69%  '
# Row 37
$ws.Range('E37').Value = '  +0.25%  '
# Row 38
$ws.Range('D38').Value = '''17.71'
$ws.Range('E38').Value = '  +0.46%  '
# Row 39
$ws.Range('D39').Value = '''1.20'
$ws.Range('E39').Value = '  +2.86%  '
# Row 40
$ws.Range('E40').Value = '  +2.46%  '
# Row 41
$ws.Range('E41').Value = '  +1.54%  '
# Row 42
$ws.Range('D42').Value = '''35.50'
$ws.Range('E42').Value = '  -0.85%  '
# Row 43
$ws.Range('E43').Value = '  +2.35%  '
# Row 44
$ws.Range('E44').Value = '  +1.99%  '
# Row 45
$ws.Range('D45').Value = '''3.36'
$ws.Range('E45').Value = '  +0.79%  '
# Row 46
$ws.Range('E46').Value = '  +4.13%  '
# Row 47
$ws.Range('D47').Value = '''4.86'
$ws.Range('E47').Value = '  +3.00%  '
# Row 48
$ws.Range('E48').Value = '  +1.27%  '
# Row 49
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '''0.552'
$ws.Range('E49').Value = '  +0.87%  '
# Row 50
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').Value = '''245.24'
$ws.Range('E50').Value = '  +2.64%  '
# Row 51
$ws.Range('D51').Value = '''0.0488'
$ws.Range('E51').Value = '  +2.53%  '
